# Update countries & provincias Spain
#
# Refresh the COVID-19 "Pais" worksheet with the latest per-country
# figures (Casos totales / Nuevos casos / Casos activos / Recuperados /
# Casos criticos / Muertes hoy / Muertes) and the "Datos actualizados"
# timestamp in A1.
#
# The sheet is kept sorted by "Casos totales" (column B) descending, so
# a handful of countries whose totals crossed a neighbour's total as a
# result of this refresh (Suazilandia / Benin / Ruanda / Jordania /
# Letonia, plus the Lesoto/Seychelles and Groenlandia/Islas Malvinas
# ties) are re-seated into their new sorted position by rewriting the
# affected rows in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1 - Datos actualizados a 9 de Julio de 2020 a las 19:33
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 9 de Julio de 2020 a las 20:50"

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 3188600
$ws.Cells.Item(4, 3).Value = 29668
$ws.Cells.Item(4, 4).Value = 1407229
$ws.Cells.Item(4, 5).Value = 1646025
$ws.Cells.Item(4, 7).Value = 484
$ws.Cells.Item(4, 8).Value = 135346

# Row 6 - India
$ws.Cells.Item(6, 2).Value = 794855
$ws.Cells.Item(6, 3).Value = 25803
$ws.Cells.Item(6, 4).Value = 495960
$ws.Cells.Item(6, 5).Value = 277272
$ws.Cells.Item(6, 7).Value = 479
$ws.Cells.Item(6, 8).Value = 21623

# Row 19 - Alemania
$ws.Cells.Item(19, 2).Value = 198952
$ws.Cells.Item(19, 3).Value = 187
$ws.Cells.Item(19, 5).Value = 6229

# Row 21 - Francia
$ws.Cells.Item(21, 2).Value = 170094
$ws.Cells.Item(21, 3).Value = 621
$ws.Cells.Item(21, 4).Value = 78170
$ws.Cells.Item(21, 5).Value = 61945
$ws.Cells.Item(21, 7).Value = 14
$ws.Cells.Item(21, 8).Value = 29979

# Row 23 - Canada
$ws.Cells.Item(23, 2).Value = 106742
$ws.Cells.Item(23, 3).Value = 308
$ws.Cells.Item(23, 5).Value = 27493

# Row 50 - Barein
$ws.Cells.Item(50, 5).Value = 4756
$ws.Cells.Item(50, 7).Value = 4
$ws.Cells.Item(50, 8).Value = 102

# Row 55 - Irlanda
$ws.Cells.Item(55, 2).Value = 25565
$ws.Cells.Item(55, 3).Value = 23
$ws.Cells.Item(55, 5).Value = 458
$ws.Cells.Item(55, 7).Value = 5
$ws.Cells.Item(55, 8).Value = 1743

# Row 62 - Argelia
$ws.Cells.Item(62, 2).Value = 17808
$ws.Cells.Item(62, 3).Value = 460
$ws.Cells.Item(62, 4).Value = 12637
$ws.Cells.Item(62, 5).Value = 4183
$ws.Cells.Item(62, 7).Value = 10
$ws.Cells.Item(62, 8).Value = 988

# Row 71 - Uzbekistan
$ws.Cells.Item(71, 5).Value = 4314
$ws.Cells.Item(71, 7).Value = 6
$ws.Cells.Item(71, 8).Value = 51

# Row 72 - Sudan
$ws.Cells.Item(72, 2).Value = 10158
$ws.Cells.Item(72, 3).Value = 74
$ws.Cells.Item(72, 5).Value = 4443
$ws.Cells.Item(72, 7).Value = 5
$ws.Cells.Item(72, 8).Value = 641

# Row 95 - Republica de Yibuti
$ws.Cells.Item(95, 2).Value = 4955
$ws.Cells.Item(95, 3).Value = 66
$ws.Cells.Item(95, 4).Value = 4671
$ws.Cells.Item(95, 5).Value = 228
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 56

# Row 109 - Maldivas
$ws.Cells.Item(109, 2).Value = 2553
$ws.Cells.Item(109, 3).Value = 36
$ws.Cells.Item(109, 4).Value = 2227
$ws.Cells.Item(109, 5).Value = 313

# Row 131 - Benin
$ws.Cells.Item(131, 1).Value = "Suazilandia"
$ws.Cells.Item(131, 2).Value = 1213
$ws.Cells.Item(131, 3).Value = 75
$ws.Cells.Item(131, 4).Value = 609
$ws.Cells.Item(131, 5).Value = 587
$ws.Cells.Item(131, 7).Value = 3
$ws.Cells.Item(131, 8).Value = 17

# Row 132 - Ruanda
$ws.Cells.Item(132, 1).Value = "Benin"
$ws.Cells.Item(132, 2).Value = 1199
$ws.Cells.Item(132, 4).Value = 333
$ws.Cells.Item(132, 5).Value = 845
$ws.Cells.Item(132, 8).Value = 21

# Row 133 - Jordania
$ws.Cells.Item(133, 1).Value = "Ruanda"
$ws.Cells.Item(133, 2).Value = 1194
$ws.Cells.Item(133, 4).Value = 610
$ws.Cells.Item(133, 5).Value = 581
$ws.Cells.Item(133, 8).Value = 3

# Row 134 - Letonia
$ws.Cells.Item(134, 1).Value = "Jordania"
$ws.Cells.Item(134, 2).Value = 1169
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 982
$ws.Cells.Item(134, 5).Value = 177
$ws.Cells.Item(134, 8).Value = 10

# Row 135 - Suazilandia
$ws.Cells.Item(135, 1).Value = "Letonia"
$ws.Cells.Item(135, 2).Value = 1154
$ws.Cells.Item(135, 3).Value = 13
$ws.Cells.Item(135, 4).Value = 1019
$ws.Cells.Item(135, 5).Value = 105
$ws.Cells.Item(135, 8).Value = 30

# Row 139 - Republica de Chipre
$ws.Cells.Item(139, 2).Value = 1010
$ws.Cells.Item(139, 3).Value = 2
$ws.Cells.Item(139, 5).Value = 152

# Row 184 - Lesoto
$ws.Cells.Item(184, 1).Value = "Seychelles"

# Row 185 - Seychelles
$ws.Cells.Item(185, 1).Value = "Lesoto"

# Row 209 - Groenlandia
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"

# Row 210 - Islas Malvinas
$ws.Cells.Item(210, 1).Value = "Groenlandia"

